$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlRight = -4152

# --- Row 3 (year headers): extend D3:J3 style (12) into K3:M3 ---
$ws.Range("J3").Copy() | Out-Null
$ws.Range("K3:M3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("K3").Value = 2020
$ws.Range("L3").Value = 2021
$ws.Range("M3").Value = 2022

# --- Row 4 (data, style 14): extend into K4:M4 ---
$ws.Range("J4").Copy() | Out-Null
$ws.Range("K4:M4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("K4").Value = 308
$ws.Range("L4").Value = 212.1
$ws.Range("M4").Value = 723.8

# --- Row 5: base style (14) extended E5:D5 into J5:M5, then values "-" in E5:L5, then right-align E5:M5 (new style 17) ---
$ws.Range("D5").Copy() | Out-Null
$ws.Range("E5:M5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E5:L5").Value = "-"
$ws.Range("E5:M5").HorizontalAlignment = $xlRight

# --- Row 6 (style 14): extend I6 format into J6:M6 ---
$ws.Range("I6").Copy() | Out-Null
$ws.Range("J6:M6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J6").Value = 9.8000000000000007
$ws.Range("K6").Value = 9.8000000000000007
$ws.Range("L6").Value = 9.8000000000000007
# M6 stays blank

# --- Row 7 (style 16): extend J7 format into K7:M7 (J7 already style 16, only needs a value) ---
$ws.Range("I7").Copy() | Out-Null
$ws.Range("K7:M7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J7").Value = 64
$ws.Range("K7").Value = 64
$ws.Range("L7").Value = 64
$ws.Range("M7").Value = 64

# --- Clear clipboard marching ants / selection ---
$ws.Range("M14").Select() | Out-Null
